$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the duplicated trailing columns X:AG (the sheet previously repeated
#    the last 10 "Pairs" columns a second time out past column W).
$ws.Range("X1:AG19").Delete()

# 2) Re-order the HKL index labels shown in row 2 (C2:M2). The underlying set of
#    labels is unchanged, only their left-to-right order changes.
$ws.Range("C2").Value = "[2, 2, 2]"
$ws.Range("D2").Value = "[1, 1, 1]"
$ws.Range("E2").Value = "[3, 1, 1]"
$ws.Range("F2").Value = "[3, 3, 1]"
$ws.Range("G2").Value = "[4, 2, 2]"
$ws.Range("H2").Value = "[5, 1, 1]"
$ws.Range("I2").Value = "[4, 0, 0]"
$ws.Range("J2").Value = "[4, 2, 0]"
$ws.Range("K2").Value = "[2, 2, 0]"
$ws.Range("L2").Value = "[2, 0, 0]"
$ws.Range("M2").Value = "[3, 3, 3]"

# 3) Insert the new "Holden" scheme rows between "OffsetATD" (row 15) and the
#    "HexGrid" rows (previously rows 16:19). Shift the HexGrid rows down to
#    20:23 and add the 4 new Holden rows at 16:19.
$ws.Rows("16:19").Insert()

$holdenNames = @("Holden2.5", "Holden5", "Holden10", "Holden15")
for ($i = 0; $i -lt 4; $i++) {
    $r = 16 + $i
    $ws.Cells.Item($r, 1).Value = 13 + $i + 1
    $ws.Cells.Item($r, 1).HorizontalAlignment = -4108
    $ws.Cells.Item($r, 1).VerticalAlignment = -4160
    $ws.Cells.Item($r, 1).Borders.LineStyle = 1
    $ws.Cells.Item($r, 1).Font.Bold = $true
    $ws.Cells.Item($r, 2).Value = $holdenNames[$i]
    for ($c = 3; $c -le 23; $c++) {
        $ws.Cells.Item($r, $c).Value = 1
    }
}

# 4) Renumber the A column for the shifted-down HexGrid rows (now 20:23) and
#    for any rows below that used to continue the sequence.
for ($r = 20; $r -le 23; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
